# CORELIMS-98 - production addition of F3
# Renumber the DRW drawer barcodes in column C from DRW562-DRW596 to DRW611-DRW645.
# Each drawer location barcode (column C) is repeated across 4 consecutive rows
# (FB1-FB4 sub-rows). Row 1 is the header row; data starts at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startOld = 562
$startNew = 611
$groupCount = 35
$firstDataRow = 2
$rowsPerGroup = 4

for ($i = 0; $i -lt $groupCount; $i++) {
    $groupStartRow = $firstDataRow + ($i * $rowsPerGroup)
    $newValue = "DRW" + ($startNew + $i)
    for ($j = 0; $j -lt $rowsPerGroup; $j++) {
        $row = $groupStartRow + $j
        $ws.Cells.Item($row, 3).Value = $newValue
    }
}
